# "fixed bag namespace issues"
#
# Backlog sheet maintenance pass:
#   - clear a stray Phase value on row 13
#   - correct Priority on rows 20/21 (3 -> 2)
#   - fill in previously-blank Priority values for rows 50-53 and 55
#   - append five new backlog items (rows 56-60) with their topic/description text
#   - re-hide every task row whose Status column is already set (complete/postponed),
#     leaving the still-open (blank Status) rows showing
#   - grow the "tasks" table / AutoFilter to cover the new rows and refresh the
#     Status-is-blank filter on the table
#   - leave the selection on the newly added block

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 13: Phase (G) value was cleared ---------------------------------
$ws.Range("G13").ClearContents()

# --- rows 20 & 21: Priority (H) 3 -> 2 ------------------------------------
$ws.Range("H20").Value = 2
$ws.Range("H21").Value = 2

# --- rows 50-53 & 55: Priority (H) filled in -------------------------------
$ws.Range("H50").Value = 3
$ws.Range("H51").Value = 4
$ws.Range("H52").Value = 4
$ws.Range("H53").Value = 2
$ws.Range("H55").Value = 4

# --- new rows 56-60: copy row 55's formatting, then fill in the content ---
$ws.Range("B55:H55").Copy()
$ws.Range("B56:H60").PasteSpecial(-4122)

$ws.Range("B56").Value = 51
$ws.Range("C56").Value = "Vaerydian"
$ws.Range("D56").Value = "UI"
$ws.Range("E56").Value = "character inventory ui"
$ws.Range("H56").Value = 3

$ws.Range("B57").Value = 52
$ws.Range("C57").Value = "Vaerydian"
$ws.Range("D57").Value = "UI"
$ws.Range("E57").Value = "character stat ui"
$ws.Range("H57").Value = 3

$ws.Range("B58").Value = 53
$ws.Range("C58").Value = "Vaerydian"
$ws.Range("D58").Value = "NPC Factory"
$ws.Range("E58").Value = "need to create more types of mobs"
$ws.Range("H58").Value = 4

$ws.Range("B59").Value = 54
$ws.Range("C59").Value = "Vaerydian"
$ws.Range("D59").Value = "Factories"
$ws.Range("E59").Value = "need to define method for factories to assemble mob templates"
$ws.Range("H59").Value = 2
$ws.Rows.Item(59).RowHeight = 30

$ws.Range("B60").Value = 55
$ws.Range("C60").Value = "Vaerydian"
$ws.Range("D60").Value = "Systems"
$ws.Range("E60").Value = "Need to define a Trigger system and Trigger Component"
$ws.Rows.Item(60).RowHeight = 30

# --- re-hide rows that already have a Status (F) value ---------------------
for ($r = 5; $r -le 55; $r++) {
    $status = $ws.Cells.Item($r, 6).Text
    if ($status -ne "") {
        $ws.Rows.Item($r).Hidden = $true
    } else {
        $ws.Rows.Item($r).Hidden = $false
    }
}

# --- grow the table to include the new rows and refresh the filter ---------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B4:H60"))
$lo.Range.AutoFilter(5, @(""), 7)

# --- leave selection on the newly-entered block -----------------------------
$ws.Range("O21").Select()
